$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.4698579619822297

$ws.Range("H4").Value = [double]"1.222175083717947e-17"
$ws.Range("I4").Value = 0.4698579619822297
$ws.Range("L4").Value = 0.12

$ws.Range("H6").Value = 0.6673333333333329
$ws.Range("I6").Value = 0.2085278213706131
$ws.Range("L6").Value = 0.1673333333333329

$ws.Range("I7").Value = 0.4519274383469447

$ws.Range("H10").Value = 0.6673333333333329
$ws.Range("I10").Value = 0.1962864880001979
$ws.Range("L10").Value = 0.1826666666666671

$ws.Range("H11").Value = 0.6673333333333329
$ws.Range("I11").Value = 0.1795051197566921
$ws.Range("L11").Value = 0.1826666666666671
